# Custom properties - bug fixes
#
# 1) "AMSIN" sheet (sprint-run history): append the two most recent sprint
#    runs (rows 68 & 69) that were missing from the exported history.
# 2) "AMS" sheet: fix row 41 ("175live") - the recorded run time had been
#    truncated/rounded; restore the precise value.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "AMSIN": add rows 68 and 69
# ---------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

# Row 68 - "176firsttrail"
# A leading apostrophe keeps the "Run Date" as literal text (e.g.
# "2023-04-18") instead of letting it be auto-parsed into a date serial.
$wsAmsin.Cells.Item(68, 1).Value = "'2023-04-18"
$wsAmsin.Cells.Item(68, 2).Value = 45034.55057184028
$wsAmsin.Cells.Item(68, 3).Value = "176firsttrail"
$wsAmsin.Cells.Item(68, 4).Value = 75
$wsAmsin.Cells.Item(68, 5).Value = 70
$wsAmsin.Cells.Item(68, 6).Value = 5
$wsAmsin.Cells.Item(68, 7).Value = 2.62

# Row 69 - "176fstscndtr"
$wsAmsin.Cells.Item(69, 1).Value = "'2023-04-18"
$wsAmsin.Cells.Item(69, 2).Value = 45034.63043397594
$wsAmsin.Cells.Item(69, 3).Value = "176fstscndtr"
$wsAmsin.Cells.Item(69, 4).Value = 75
$wsAmsin.Cells.Item(69, 5).Value = 71
$wsAmsin.Cells.Item(69, 6).Value = 4
$wsAmsin.Cells.Item(69, 7).Value = 2.53

# Column B ("Run Time") keeps the same date/time display as the rest of
# the column - copy that formatting down onto the two new rows.
$wsAmsin.Range("B67").Copy()
$wsAmsin.Range("B68:B69").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Sheet "AMS": fix row 41 ("175live")
# ---------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

# Correct the recorded run-time precision.
$wsAms.Cells.Item(41, 2).Value = 45016.79465128472
